# Insert 9 new daily rows (2019-11-18 .. 2019-11-28) into the AEMULUS (0181)
# price history sheet, just before the existing 2019-11-29 row, shifting all
# subsequent rows down by 9 (old row 1019 -> new row 1028, ..., old row 1091
# -> new row 1100). The sheet dimension grows from I1091 to I1100.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 1019-1091 down by inserting 9 blank rows at row 1019.
$ws.Rows("1019:1027").Insert()

# Data for the newly inserted rows: timestamp, date, id, name, open, high, low, close, vol
$newRows = @(
    @(1019, 1574035200, "2019-11-18", "0181", "AEMULUS", 0.265, 0.27,  0.265, 0.265, 584700),
    @(1020, 1574121600, "2019-11-19", "0181", "AEMULUS", 0.265, 0.265, 0.26,  0.265, 642100),
    @(1021, 1574208000, "2019-11-20", "0181", "AEMULUS", 0.265, 0.265, 0.255, 0.26,  1616600),
    @(1022, 1574294400, "2019-11-21", "0181", "AEMULUS", 0.255, 0.26,  0.25,  0.26,  1051600),
    @(1023, 1574380800, "2019-11-22", "0181", "AEMULUS", 0.24,  0.245, 0.235, 0.24,  3724300),
    @(1024, 1574640000, "2019-11-25", "0181", "AEMULUS", 0.24,  0.24,  0.23,  0.23,  1538000),
    @(1025, 1574726400, "2019-11-26", "0181", "AEMULUS", 0.23,  0.23,  0.225, 0.225, 835700),
    @(1026, 1574812800, "2019-11-27", "0181", "AEMULUS", 0.225, 0.225, 0.22,  0.225, 1250700),
    @(1027, 1574899200, "2019-11-28", "0181", "AEMULUS", 0.225, 0.225, 0.22,  0.22,  449400)
)

foreach ($r in $newRows) {
    $row = $r[0]

    # Columns B (date) and C (id) are stored as literal text in this
    # workbook ("2019-11-18", "0181") -- force text format first so Excel
    # doesn't auto-coerce them into a date serial / number (losing the
    # leading zero on the id).
    $ws.Cells.Item($row, 2).NumberFormat = "@"
    $ws.Cells.Item($row, 3).NumberFormat = "@"

    $ws.Cells.Item($row, 1).Value = $r[1]
    $ws.Cells.Item($row, 2).Value = $r[2]
    $ws.Cells.Item($row, 3).Value = $r[3]
    $ws.Cells.Item($row, 4).Value = $r[4]
    $ws.Cells.Item($row, 5).Value = $r[5]
    $ws.Cells.Item($row, 6).Value = $r[6]
    $ws.Cells.Item($row, 7).Value = $r[7]
    $ws.Cells.Item($row, 8).Value = $r[8]
    $ws.Cells.Item($row, 9).Value = $r[9]
}
